# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets,
# matching the refreshed output data captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 55
$ws1.Range("F4").Value = 4641
$ws1.Range("F5").Value = 1843
$ws1.Range("F6").Value = 135
$ws1.Range("F8").Value = 3112
$ws1.Range("F10").Value = 587
$ws1.Range("F11").Value = 266
$ws1.Range("F12").Value = 630
$ws1.Range("F13").Value = 537
$ws1.Range("F14").Value = 529
$ws1.Range("F17").Value = 1779
$ws1.Range("F18").Value = 1335
$ws1.Range("F20").Value = 1608
$ws1.Range("F21").Value = 6
$ws1.Range("F24").Value = 8
$ws1.Range("F26").Value = 535
$ws1.Range("F33").Value = 3831
$ws1.Range("F34").Value = 763
$ws1.Range("F35").Value = 73
$ws1.Range("F36").Value = 836
$ws1.Range("F38").Value = 1826

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 55
$ws4.Range("F4").Value = 4641
$ws4.Range("F5").Value = 1843
$ws4.Range("F6").Value = 135
$ws4.Range("F8").Value = 3112
$ws4.Range("F10").Value = 587
$ws4.Range("F11").Value = 266
$ws4.Range("F12").Value = 630
$ws4.Range("F13").Value = 537
$ws4.Range("F14").Value = 529
$ws4.Range("F18").Value = 1779
$ws4.Range("F19").Value = 1335
$ws4.Range("F21").Value = 1608
$ws4.Range("F22").Value = 6
$ws4.Range("F25").Value = 8
$ws4.Range("F27").Value = 535
$ws4.Range("F34").Value = 3831
$ws4.Range("F36").Value = 763
$ws4.Range("F37").Value = 73
$ws4.Range("F38").Value = 836
$ws4.Range("F40").Value = 1826
